$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header labels to include a space
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"

# Update the active selection
$ws.Range("D4").Select()
